# Adds Multilayer Perceptron (Random search) results to rows 22/23 of all
# 8 worksheets (Single Model block, column groups B:H and J:P), matching
# the commit "add MLP and changed read_train".
$wb = $excel.ActiveWorkbook

# Sheet 1: Option 1 - LR1 - DN1 (70-30)
$ws = $wb.Worksheets.Item(1)
$ws.Range("B22").Value = 'Random'
$ws.Range("C22").Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range("D22").Value = 0.074021115899086
$ws.Range("E22").Value = 0.2039611786603928
$ws.Range("F22").Value = 2.475558757781982
$ws.Range("G22").Value = 0.2720682192007843
$ws.Range("H22").Value = 40.00687301158905
$ws.Range("J22").Value = 'Random'
$ws.Range("K22").Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range("L22").Value = 0.07532788068056107
$ws.Range("M22").Value = 0.208153635263443
$ws.Range("N22").Value = 2.494367122650146
$ws.Range("O22").Value = 0.2744592514027557
$ws.Range("P22").Value = 40.09046256542206

# Sheet 2: Option 1 - LR1 - DN2 (70-30)
$ws = $wb.Worksheets.Item(2)
$ws.Range("B23").Value = 'Random'
$ws.Range("C23").Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range("D23").Value = 0.9456392526626587
$ws.Range("E23").Value = 0.5076738595962524
$ws.Range("F23").Value = 0.505184531211853
$ws.Range("G23").Value = 0.9724398452668723
$ws.Range("H23").Value = 24.99835640192032
$ws.Range("J23").Value = 'Random'
$ws.Range("K23").Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range("L23").Value = 0.9287254214286804
$ws.Range("M23").Value = 0.5096931457519531
$ws.Range("N23").Value = 0.5012574791908264
$ws.Range("O23").Value = 0.9637040113171058
$ws.Range("P23").Value = 24.97094869613647

# Sheet 3: Option 1 - LR2 - DN1 (70-30)
$ws = $wb.Worksheets.Item(3)
$ws.Range("B22").Value = 'Random'
$ws.Range("C22").Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.02}'
$ws.Range("D22").Value = 55.48139190673828
$ws.Range("E22").Value = 5.74082612991333
$ws.Range("F22").Value = 2.712109327316284
$ws.Range("G22").Value = 7.44858321472871
$ws.Range("H22").Value = 43.76845955848694
$ws.Range("J22").Value = 'Random'
$ws.Range("K22").Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.02}'
$ws.Range("L22").Value = 54.68329620361328
$ws.Range("M22").Value = 5.723203182220459
$ws.Range("N22").Value = 2.702246189117432
$ws.Range("O22").Value = 7.39481549490001
$ws.Range("P22").Value = 43.56915950775146

# Sheet 4: Option 1 - LR2 - DN2 (70-30)
$ws = $wb.Worksheets.Item(4)
$ws.Range("B22").Value = 'Random'
$ws.Range("C22").Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.05}'
$ws.Range("D22").Value = 62.94614791870117
$ws.Range("E22").Value = 6.126215934753418
$ws.Range("F22").Value = 2.863848686218262
$ws.Range("G22").Value = 7.933860845685484
$ws.Range("H22").Value = 45.27258276939392
$ws.Range("J22").Value = 'Random'
$ws.Range("K22").Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.05}'
$ws.Range("L22").Value = 37.5901985168457
$ws.Range("M22").Value = 4.847194671630859
$ws.Range("N22").Value = 1.935303092002869
$ws.Range("O22").Value = 6.131084611783278
$ws.Range("P22").Value = 41.4331465959549

# Sheet 5: Option 1 - NLR1 - DN1 (70-30)
$ws = $wb.Worksheets.Item(5)
$ws.Range("B22").Value = 'Random'
$ws.Range("C22").Value = '{''module__num_units'': 50, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range("D22").Value = 0.08401425182819366
$ws.Range("E22").Value = 0.2262209206819534
$ws.Range("F22").Value = 29832033337344
$ws.Range("G22").Value = 0.2898521206204875
$ws.Range("H22").Value = 43.09643507003784
$ws.Range("J22").Value = 'Random'
$ws.Range("K22").Value = '{''module__num_units'': 20, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range("L22").Value = 0.08294574171304703
$ws.Range("M22").Value = 0.2236228734254837
$ws.Range("N22").Value = 29937593483264
$ws.Range("O22").Value = 0.2880030237914995
$ws.Range("P22").Value = 43.08747351169586

# Sheet 6: Option 1 - NLR1 - DN2 (70-30)
$ws = $wb.Worksheets.Item(6)
$ws.Range("B22").Value = 'Random'
$ws.Range("C22").Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range("D22").Value = 1.20376193523407
$ws.Range("E22").Value = 0.6018911004066467
$ws.Range("F22").Value = 0.6955813765525818
$ws.Range("G22").Value = 1.097160852033133
$ws.Range("H22").Value = 27.0847350358963
$ws.Range("J22").Value = 'Random'
$ws.Range("K22").Value = '{''module__num_units'': 20, ''module__activation_func'': Tanh(), ''lr'': 0.1}'
$ws.Range("L22").Value = 1.130769371986389
$ws.Range("M22").Value = 0.5777952075004578
$ws.Range("N22").Value = 0.7077444791793823
$ws.Range("O22").Value = 1.063376401838215
$ws.Range("P22").Value = 26.84087157249451

# Sheet 7: Option 1 - NLR2 - DN1 (70-30)
$ws = $wb.Worksheets.Item(7)
$ws.Range("B22").Value = 'Random'
$ws.Range("C22").Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range("D22").Value = 0.08315457403659821
$ws.Range("E22").Value = 0.2205072939395905
$ws.Range("F22").Value = 2.589795589447021
$ws.Range("G22").Value = 0.2883653481897542
$ws.Range("H22").Value = 40.86825549602509

# Sheet 8: Option 1 - NLR2 - DN2 (70-30)
$ws = $wb.Worksheets.Item(8)
$ws.Range("B22").Value = 'Random'
$ws.Range("C22").Value = '{''module__num_units'': 10, ''module__activation_func'': ReLU(), ''lr'': 0.1}'
$ws.Range("D22").Value = 0.9467169046401978
$ws.Range("E22").Value = 0.5165208578109741
$ws.Range("F22").Value = 0.5075122714042664
$ws.Range("G22").Value = 0.9729937844817909
$ws.Range("H22").Value = 25.01609921455383
